$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-03-13 Wednesday" "2024-03-14 Thursday"

Replace-Text "89×18=1602" "75×21=1575"
Replace-Text "33×89=2937" "40×23=920"
Replace-Text "41×65=2665" "30×11=330"
Replace-Text "50×80=4000" "25×29=725"
Replace-Text "75×49=3675" "38×30=1140"

Replace-Text "99×80=7920" "26×43=1118"
Replace-Text "57×16=912" "91×92=8372"
Replace-Text "88×20=1760" "80×66=5280"
Replace-Text "41×25=1025" "92×51=4692"
Replace-Text "16×61=976" "56×95=5320"

Replace-Text "80×87=6960" "60×80=4800"
Replace-Text "29×45=1305" "29×56=1624"
Replace-Text "83×48=3984" "85×37=3145"
Replace-Text "20×85=1700" "22×87=1914"
Replace-Text "76×92=6992" "91×62=5642"

Replace-Text "22×45=990" "79×26=2054"
Replace-Text "62×69=4278" "22×55=1210"
Replace-Text "18×25=450" "95×80=7600"
Replace-Text "73×46=3358" "53×48=2544"
Replace-Text "35×65=2275" "52×94=4888"

Replace-Text "18×52=936" "35×97=3395"
Replace-Text "70×24=1680" "24×69=1656"
Replace-Text "32×94=3008" "13×75=975"
Replace-Text "42×69=2898" "13×28=364"
Replace-Text "63×58=3654" "39×69=2691"
